# SW620_noCTRL_meas.xlsx bug fix:
#  - Sheet1 had stray leftover rows (45:87) with only column A populated
#    (an artifact of a previous paste); remove them so the sheet's used
#    range matches the real data (A1:N44), and leave the view focused on
#    Sheet1 instead of Sheet3.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Drop the leftover junk rows below the real data table.
$ws1.Rows("45:87").Delete()

# Sheet1 becomes the active/visible sheet (was Sheet3).
$ws1.Activate()
$ws1.Range("G58").Select()
